# Update cryptos list with latest scraped prices/volumes (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds plain scraped text (e.g. "4.17", "26.603.05") that
# must stay text, not get reinterpreted as numbers. Values with two dots
# (thousands+decimal, e.g. "26.603.05") never parse as a number so they are
# safe to assign directly. Plain-decimal-looking values (e.g. "4.17") would
# otherwise be auto-converted to a number by Excel, so those cells are forced
# to Text format first, exactly as Excel would require for literal entry.
$textPriceCells = "D5","D10","D14","D15","D16","D19","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D38","D39","D45","D47","D50","D51"
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.603.05"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.639.41"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.17%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.94"
$ws.Range("E5").Value = "  +0.36%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.15%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.71%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.14%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.02"
$ws.Range("E10").Value = "  -0.10%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.09%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.865.03"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.646.53"
$ws.Range("E13").Value = "  +0.54%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "4.17"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.529"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "65.50"
$ws.Range("E16").Value = "  +3.67%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.614.46"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +0.39%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "217.45"
$ws.Range("E19").Value = "  +3.62%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.17%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "4.34"
$ws.Range("E21").Value = "  +0.74%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "6.30"
$ws.Range("E22").Value = "  +1.43%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.40"
$ws.Range("E23").Value = "  -0.66%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.17"
$ws.Range("E24").Value = "  +13.39%  "

# Row 25 - Monero
$ws.Range("D25").Value = "147.55"

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.21%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.41%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "6.93"
$ws.Range("E28").Value = "  +0.45%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.70"
$ws.Range("E29").Value = "  +1.91%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  -0.83%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -0.68%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +3.85%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +1.70%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.271.90"
$ws.Range("E34").Value = "  +8.71%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +0.61%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +1.18%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +2.33%  "

# Row 38 & 39 - ImmutableX and ARBITRUM swapped ranking order
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.805"
$ws.Range("E38").Value = "  -0.81%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.513"
$ws.Range("E39").Value = "  +1.51%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.16%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -1.59%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  +0.96%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -0.30%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.774.71"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45 - Quant
$ws.Range("D45").Value = "93.71"
$ws.Range("E45").Value = "  +1.34%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +3.10%  "

# Row 47 - Aave
$ws.Range("D47").Value = "55.34"
$ws.Range("E47").Value = "  +1.16%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  -1.92%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.29%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "7.63"
$ws.Range("E50").Value = "  +0.60%  "

# Row 51 - Algorand
$ws.Range("D51").Value = "0.0964"
$ws.Range("E51").Value = "  +2.43%  "
